$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) values between row 2<->4 and row 3<->5
$ws.Range("D2").Value = 44305
$ws.Range("D3").Value = 44305
$ws.Range("D4").Value = 44309
$ws.Range("D5").Value = 44309

# Swap the "Volumen" (M) values between row 2<->4 and row 3<->5
$ws.Range("M2").Value = 50
$ws.Range("M3").Value = 60
$ws.Range("M4").Value = 40
$ws.Range("M5").Value = 70
